$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new data row (row 41) following the existing date/mileage series
$ws.Range("A41").NumberFormat = "[$-40C]d\-mmm\-yy;@"
$ws.Range("A41").Value = 43815

$ws.Range("B41").Value = 2393

# Update the active selection to reflect the next empty row, as Excel does
# after data entry
$ws.Application.Goto($ws.Range("B42"))
